$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B cells (rows 2-38) ---
$ws.Range("B2").Value = "NSE:ARE&M"
$ws.Range("B3").Value = "NSE:AVTNPL"
$ws.Range("B4").Value = "NSE:BAJAJHIND"
$ws.Range("B5").Value = "NSE:BANKBEES"
$ws.Range("B6").Value = "NSE:BEARDSELL"
$ws.Range("B7").Value = "NSE:BHARATGEAR"
$ws.Range("B8").Value = "NSE:CINEVISTA"
$ws.Range("B9").Value = "NSE:CLSEL"
$ws.Range("B10").Value = "NSE:COFORGE"
$ws.Range("B11").Value = "NSE:CONSUMBEES"
$ws.Range("B12").Value = "NSE:DBCORP"
$ws.Range("B13").Value = "NSE:DELTAMAGNT"
$ws.Range("B14").Value = "NSE:DEN"
$ws.Range("B15").Value = "NSE:EBBETF0425"
$ws.Range("B16").Value = "NSE:ELDEHSG"
$ws.Range("B17").Value = "NSE:ESTER"
$ws.Range("B18").Value = "NSE:EVEREADY"
$ws.Range("B19").Value = "NSE:FIBERWEB"
$ws.Range("B20").Value = "NSE:FORTIS"
$ws.Range("B21").Value = "NSE:GAEL"
$ws.Range("B22").Value = "NSE:GHCL"
$ws.Range("B23").Value = "NSE:GHCLTEXTIL"
$ws.Range("B24").Value = "NSE:GLAND"
$ws.Range("B25").Value = "NSE:HAL"
$ws.Range("B26").Value = "NSE:HDFCBSE500"
$ws.Range("B27").Value = "NSE:HDFCNIFBAN"
$ws.Range("B28").Value = "NSE:HEIDELBERG"
$ws.Range("B29").Value = "NSE:HERANBA"
$ws.Range("B30").Value = "NSE:HEXATRADEX"
$ws.Range("B31").Value = "NSE:HIL"
$ws.Range("B32").Value = "NSE:HONDAPOWER"
$ws.Range("B33").Value = "NSE:IOLCP"
$ws.Range("B34").Value = "NSE:IONEXCHANG"
$ws.Range("B35").Value = "NSE:J&KBANK"
$ws.Range("B36").Value = "NSE:KALAMANDIR"
$ws.Range("B37").Value = "NSE:KEC"
$ws.Range("B38").Value = "NSE:KIMS"

# --- Update column C cells ---
$ws.Range("C2").Value = "NSE:ARVSMART"
$ws.Range("C3").Value = "NSE:OAL"

# --- Clear column C4 (value removed) ---
$ws.Range("C4").Value = ""

# --- Add new rows 39-56, copying formatting from row 38 (style s="1" on col A) ---
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "NSE:LAGNAM"
$ws.Range("A38").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "NSE:MASTEK"
$ws.Range("A38").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "NSE:MEDICAMEQ"
$ws.Range("A38").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "NSE:MHRIL"
$ws.Range("A38").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "NSE:MOQUALITY"
$ws.Range("A38").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "NSE:NILAINFRA"
$ws.Range("A38").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "NSE:NOCIL"
$ws.Range("A38").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "NSE:NUVOCO"
$ws.Range("A38").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = "NSE:ORIENTLTD"
$ws.Range("A38").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "NSE:OSWALGREEN"
$ws.Range("A38").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = "NSE:PALASHSECU"
$ws.Range("A38").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "NSE:PFIZER"
$ws.Range("A38").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "NSE:PKTEA"
$ws.Range("A38").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = "NSE:PRECOT"
$ws.Range("A38").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = "NSE:PRICOLLTD"
$ws.Range("A38").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = "NSE:PURVA"
$ws.Range("A38").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = "NSE:REPL"
$ws.Range("A38").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = "NSE:SAIL"

$excel.CutCopyMode = 0
